$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

# Update the password value from "Admin123" to "admin123"
$ws.Range("B2").Value = "admin123"

# Update the active selection on the sheet to C2
$ws.Range("C2").Select()
